{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 0: date/title line (two runs separated by a line break) ---\nparas.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 11.03.25\\u000bBeyond Matryoshka: Revisiting Sparse Coding for Adaptive Representation\", \"Replace\");\n\n// --- Paragraphs 1-7: full text replacement ---\nparas.items[1].insertText(\"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e7\u05e6\u05e8\u05d4 \u05e9\u05dc \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05db\u05dc\u05d9\u05dc \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d4\u05e4\u05e7\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d1\u05de\u05d9\u05de\u05d3 \u05e0\u05de\u05d5\u05da \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea Matryoshka embeddings. \u05de\u05d4 \u05de\u05d9\u05d5\u05d7\u05d3 \u05d1\u05e9\u05d9\u05d8\u05d4 \u05d6\u05d5 - \u05d4\u05d9\u05d0 \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d6\u05d4 \u05d1\u05db\u05de\u05d4 \u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05de\u05db\u05de\u05d4 \u05d2\u05d3\u05dc\u05d9\u05dd (\u05e0\u05d2\u05d9\u05d3 8\u05bf, 16, 32, 64 \u05d5-128) \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d6\u05de\u05df. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05e0\u05d9\u05d7\u05d4 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05de\u05ea\u05d5\u05d9\u05d2 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea (x, y) \u05db\u05d0\u05e9\u05e8 x \u05d4\u05d5\u05d0 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d5-y \u05d4\u05d5\u05d0 \u05d4\u05ea\u05d9\u05d5\u05d2 \u05e9\u05dc\u05d5. \", \"Replace\");\nparas.items[2].insertText(\"\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e8\u05e9\u05ea \u05e2\u05de\u05d5\u05e7\u05d4 \u05e2\u05dd \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 (\u05e8\u05d0\u05e9) \u05d4\u05de\u05de\u05e4\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05dc\u05ea\u05d9\u05d5\u05d2 \u05e9\u05dc\u05d5. \u05de\u05d4 \u05d4\u05de\u05d9\u05d5\u05d7\u05d3 \u05d1\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d4\u05d5\u05d0 \u05e9\u05d4\u05d9\u05d0 \u05de\u05d0\u05de\u05e0\u05ea \u05d1\u05d5-\u05d6\u05de\u05e0\u05d9\u05ea \u05db\u05de\u05d4 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05de\u05d9\u05e4\u05d5\u05d9 (\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5) \u05dc\u05de\u05d7\u05e8\u05d1 \u05d4\u05ea\u05d9\u05d5\u05d2 \u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05de\u05d9\u05e4\u05d5\u05d9 \u05dc\u05d5\u05e7\u05d7 m_i \u05d4\u05d0\u05d9\u05d1\u05e8\u05d9\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05de\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2(\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05d1\u05d3\u05d5\u05d2\u05de\u05d0 \u05e9\u05e0\u05ea\u05ea\u05d9 \u05e7\u05d5\u05d3\u05dd \u05de\u05d0\u05de\u05df \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05de\u05d9\u05e4\u05d5\u05d9 \u05d1\u05d2\u05d3\u05dc\u05d9\u05dd 8, 16, 32 \u05d5-64. \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1 \u05d4\u05d9\u05e0\u05d4 \u05e1\u05db\u05d5\u05dd \u05e9\u05dc \u05d4\u05dc\u05d5\u05e1\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 - \u05db\u05dc\u05d5\u05de\u05e8 \u05e0\u05d5\u05e1\u05e3 \u05dc\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd 4 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d1\u05d2\u05d3\u05dc\u05d9\u05dd 8, 16, 32 \u05d5- 64. \", \"Replace\");\nparas.items[3].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05db\u05dc\u05d9\u05dc \u05d0\u05ea \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea \u05d4\u05d6\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d7\u05dc\u05e4\u05ea\u05d4 \u05d1\u05e9\u05e0\u05d9 \u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd(\u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1 \u05dc\u05de\u05e2\u05e9\u05d1\u05d4). \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 sparse autoencoder \u05d0\u05d5 SE \u05e9\u05d1\u05de\u05e7\u05d5\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05de\u05e4\u05d5\u05ea \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d3\u05d0\u05d8\u05d4, \u05d4\u05de\u05d5\u05e4\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc, \u05dc\u05de\u05e8\u05d7\u05d1 \u05d1\u05e2\u05dc \u05de\u05d9\u05de\u05d3 \u05de\u05d0\u05d5\u05d3 \u05d2\u05d1\u05d5\u05d4 \u05d0\u05d1\u05dc \u05de\u05d0\u05d5\u05d3 \u05d3\u05dc\u05d9\u05dc \u05d5\u05d0\u05d6 \u05dc\u05d4\u05d7\u05d6\u05d9\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9. \u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05dc\u05d0 \u05de\u05d0\u05d5\u05de\u05df \u05db\u05d0\u05df \u05d0\u05dc\u05d0 \u05e8\u05e7 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 (\u05e9\u05dc SE). \u05d4\u05d0\u05dc\u05de\u05e0\u05d8 \u05d4\u05e9\u05e0\u05d9 \u05e9\u05de\u05ea\u05d5\u05d5\u05e1\u05e3 \u05e9\u05d4\u05dc\u05d5\u05e1 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9 \u05e9\u05d1\u05d0 \u05dc\u05d4\u05e8\u05d7\u05d9\u05e7 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d3\u05d0\u05d8\u05d4 \u05de\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e8\u05d7\u05d5\u05e7 \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd \u05d5\u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05de\u05d0\u05d5\u05ea\u05d4 \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4.\", \"Replace\");\nparas.items[4].insertText(\"\u05bf\u05d0\u05d6 \u05de\u05d4 \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc SE \u05db\u05d0\u05df? \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05d4\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05e9\u05de\u05d0\u05de\u05e0\u05ea \u05d0\u05ea \u05d4\u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05db\u05d0\u05df \u05d0\u05e0\u05d5 \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd top-k \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d3\u05d5\u05e8. \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05e9\u05d7\u05d6\u05e8 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9 \u05e8\u05e7 \u05e2\u05dd top-k \u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8. \u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d9\u05d3\u05d5\u05e2 \u05e2\u05dd SE \u05d4\u05d9\u05d0 \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05ea\u05d9\u05dd - \u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05e2\u05e8\u05db\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05dc\u05db\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4. \", \"Replace\");\nparas.items[5].insertText(\"\u05db\u05d3\u05d9 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05d1\u05e2\u05d9\u05d4 \u05d6\u05d5 \u05d4\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05e0\u05d9 \u05d3\u05d1\u05e8\u05d9\u05dd. \u05d4\u05d3\u05d1\u05e8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 \u05d4\u05d5\u05e1\u05e4\u05d4 \u05dc\u05d5\u05e1\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05de\u05d4 \u05e2\u05e8\u05db\u05d9\u05dd \u05e9\u05dc k \u05dc- top-k \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 (\u05d1\u05de\u05e7\u05d5\u05e8 \u05d9\u05e9 \u05e2\u05e8\u05da k \u05d0\u05d7\u05d3). \u05db\u05db\u05d4 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05d1\u05db\u05de\u05d4 \u05d2\u05d3\u05dc\u05d9\u05dd \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 (\u05d7\u05d5\u05e5 \u05de\u05d6\u05d4 \u05d0\u05d9\u05df \u05d4\u05e8\u05d1\u05d4 \u05d3\u05de\u05d9\u05d5\u05df \u05db\u05d9 \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d4\u05e4\u05d9\u05e7 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05d3\u05dc\u05d9\u05dc). \u05d4\u05d3\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e9\u05dc \u05d0\u05d9\u05d1\u05e8 \u05d4\u05de\u05e0\u05e1\u05d4 \u05dc\u05d2\u05e8\u05d5\u05dd \u05dc\u05e9\u05d2\u05d9\u05d0\u05ea \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e2\u05d1\u05d5\u05e8 top-k \u05e9\u05dc \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05de\u05ea\u05d9\u05dd (\u05e2\u05e8\u05db\u05d9\u05dd \u05d4\u05db\u05d9 \u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8) \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e9\u05d2\u05d9\u05d0\u05ea \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e9\u05dc \u05d4-top-k \u05e9\u05dc \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d5\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8. \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05e8\u05d3\u05ea \u05dc\u05e2\u05d5\u05de\u05e7 \u05d3\u05e2\u05ea\u05dd \u05dc\u05de\u05d4 \u05d6\u05d4 \u05e2\u05d5\u05d6\u05e8.\", \"Replace\");\nparas.items[6].insertText(\"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05d0\u05d9\u05d1\u05e8 \u05e9\u05dc \u05d4\u05dc\u05d5\u05e1 \u05d4\u05e0\u05d9\u05d2\u05d5\u05d3\u05d9 \u05dc\u05d6\u05d4 \u05e9\u05de\u05ea\u05d5\u05d0\u05e8 \u05d1\u05e4\u05e1\u05e7\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea\u2026.\", \"Replace\");\nparas.items[7].insertText(\"\u05d8\u05d5\u05d1, \u05e0\u05db\u05d5\u05df \u05e9\u05d4\u05d5\u05e4\u05d9\u05e2\u05d4 \u05dc\u05e0\u05d5 \u05d4\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d1\u05e9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05d1\u05d9\u05e0\u05d5 \u05dc\u05d1\u05d9\u05df \u05d4\u05de\u05d8\u05e8\u05d9\u05e9\u05e7\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05d3\u05d9 \u05e8\u05d5\u05e4\u05e3. \u05d0\u05d1\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d7\u05d5\u05e5 \u05de\u05d6\u05d4\u2026.\", \"Replace\");\nawait context.sync();\n\n// --- Paragraphs 8-22: delete entirely (removed body of the old review) ---\nfor (let i = 22; i >= 8; i--) {\n  paras.items[i].delete();\n}\nawait context.sync();\n\n// --- Paragraph 23 (now last paragraph, the link): replace URL text ---\nparas.items[23].insertText(\"https://arxiv.org/pdf/2503.01776\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Paragraph 1 (1-based): date/title line (two runs separated by a line break) ---\n$d.Paragraphs.Item(1).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 11.03.25\" + [char]11 + \"Beyond Matryoshka: Revisiting Sparse Coding for Adaptive Representation\"\n\n# --- Paragraphs 2-8 (1-based): full text replacement ---\n$d.Paragraphs.Item(2).Range.Text = \"\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e7\u05e6\u05e8\u05d4 \u05e9\u05dc \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05db\u05dc\u05d9\u05dc \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d4\u05e4\u05e7\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d1\u05de\u05d9\u05de\u05d3 \u05e0\u05de\u05d5\u05da \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea Matryoshka embeddings. \u05de\u05d4 \u05de\u05d9\u05d5\u05d7\u05d3 \u05d1\u05e9\u05d9\u05d8\u05d4 \u05d6\u05d5 - \u05d4\u05d9\u05d0 \u05de\u05d0\u05e4\u05e9\u05e8\u05ea \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d6\u05d4 \u05d1\u05db\u05de\u05d4 \u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05de\u05db\u05de\u05d4 \u05d2\u05d3\u05dc\u05d9\u05dd (\u05e0\u05d2\u05d9\u05d3 8\u05bf, 16, 32, 64 \u05d5-128) \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05d1\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d6\u05de\u05df. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05e0\u05d9\u05d7\u05d4 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05de\u05ea\u05d5\u05d9\u05d2 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea (x, y) \u05db\u05d0\u05e9\u05e8 x \u05d4\u05d5\u05d0 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d5-y \u05d4\u05d5\u05d0 \u05d4\u05ea\u05d9\u05d5\u05d2 \u05e9\u05dc\u05d5. \"\n$d.Paragraphs.Item(3).Range.Text = \"\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e8\u05e9\u05ea \u05e2\u05de\u05d5\u05e7\u05d4 \u05e2\u05dd \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 (\u05e8\u05d0\u05e9) \u05d4\u05de\u05de\u05e4\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05dc\u05ea\u05d9\u05d5\u05d2 \u05e9\u05dc\u05d5. \u05de\u05d4 \u05d4\u05de\u05d9\u05d5\u05d7\u05d3 \u05d1\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d4\u05d5\u05d0 \u05e9\u05d4\u05d9\u05d0 \u05de\u05d0\u05de\u05e0\u05ea \u05d1\u05d5-\u05d6\u05de\u05e0\u05d9\u05ea \u05db\u05de\u05d4 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05de\u05d9\u05e4\u05d5\u05d9 (\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5) \u05dc\u05de\u05d7\u05e8\u05d1 \u05d4\u05ea\u05d9\u05d5\u05d2 \u05db\u05d0\u05e9\u05e8 \u05db\u05dc \u05de\u05d9\u05e4\u05d5\u05d9 \u05dc\u05d5\u05e7\u05d7 m_i \u05d4\u05d0\u05d9\u05d1\u05e8\u05d9\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05de\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2(\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05d1\u05d3\u05d5\u05d2\u05de\u05d0 \u05e9\u05e0\u05ea\u05ea\u05d9 \u05e7\u05d5\u05d3\u05dd \u05de\u05d0\u05de\u05df \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05de\u05d9\u05e4\u05d5\u05d9 \u05d1\u05d2\u05d3\u05dc\u05d9\u05dd 8, 16, 32 \u05d5-64. \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1 \u05d4\u05d9\u05e0\u05d4 \u05e1\u05db\u05d5\u05dd \u05e9\u05dc \u05d4\u05dc\u05d5\u05e1\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 - \u05db\u05dc\u05d5\u05de\u05e8 \u05e0\u05d5\u05e1\u05e3 \u05dc\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd 4 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d1\u05d2\u05d3\u05dc\u05d9\u05dd 8, 16, 32 \u05d5- 64. \"\n$d.Paragraphs.Item(4).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05db\u05dc\u05d9\u05dc \u05d0\u05ea \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea \u05d4\u05d6\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d7\u05dc\u05e4\u05ea\u05d4 \u05d1\u05e9\u05e0\u05d9 \u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd(\u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1 \u05dc\u05de\u05e2\u05e9\u05d1\u05d4). \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 sparse autoencoder \u05d0\u05d5 SE \u05e9\u05d1\u05de\u05e7\u05d5\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05de\u05e4\u05d5\u05ea \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d3\u05d0\u05d8\u05d4, \u05d4\u05de\u05d5\u05e4\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc, \u05dc\u05de\u05e8\u05d7\u05d1 \u05d1\u05e2\u05dc \u05de\u05d9\u05de\u05d3 \u05de\u05d0\u05d5\u05d3 \u05d2\u05d1\u05d5\u05d4 \u05d0\u05d1\u05dc \u05de\u05d0\u05d5\u05d3 \u05d3\u05dc\u05d9\u05dc \u05d5\u05d0\u05d6 \u05dc\u05d4\u05d7\u05d6\u05d9\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05de\u05e8\u05d7\u05d1 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9. \u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05e6\u05de\u05d5 \u05dc\u05d0 \u05de\u05d0\u05d5\u05de\u05df \u05db\u05d0\u05df \u05d0\u05dc\u05d0 \u05e8\u05e7 \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 (\u05e9\u05dc SE). \u05d4\u05d0\u05dc\u05de\u05e0\u05d8 \u05d4\u05e9\u05e0\u05d9 \u05e9\u05de\u05ea\u05d5\u05d5\u05e1\u05e3 \u05e9\u05d4\u05dc\u05d5\u05e1 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9 \u05e9\u05d1\u05d0 \u05dc\u05d4\u05e8\u05d7\u05d9\u05e7 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d3\u05d0\u05d8\u05d4 \u05de\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e8\u05d7\u05d5\u05e7 \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd \u05d5\u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05de\u05d0\u05d5\u05ea\u05d4 \u05d4\u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d4.\"\n$d.Paragraphs.Item(5).Range.Text = \"\u05bf\u05d0\u05d6 \u05de\u05d4 \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc SE \u05db\u05d0\u05df? \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05d4\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05e9\u05de\u05d0\u05de\u05e0\u05ea \u05d0\u05ea \u05d4\u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d9\u05dd \u05db\u05d0\u05df \u05d0\u05e0\u05d5 \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd top-k \u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d3\u05d5\u05e8. \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05e9\u05d7\u05d6\u05e8 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9 \u05e8\u05e7 \u05e2\u05dd top-k \u05d0\u05dc\u05de\u05e0\u05d8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8. \u05d4\u05d1\u05e2\u05d9\u05d4 \u05d4\u05d9\u05d3\u05d5\u05e2 \u05e2\u05dd SE \u05d4\u05d9\u05d0 \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05ea\u05d9\u05dd - \u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05e2\u05e8\u05db\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05dc\u05db\u05dc \u05e4\u05d9\u05e1\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4. \"\n$d.Paragraphs.Item(6).Range.Text = \"\u05db\u05d3\u05d9 \u05dc\u05d4\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05d1\u05e2\u05d9\u05d4 \u05d6\u05d5 \u05d4\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05e0\u05d9 \u05d3\u05d1\u05e8\u05d9\u05dd. \u05d4\u05d3\u05d1\u05e8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05d5\u05d0 \u05d4\u05d5\u05e1\u05e4\u05d4 \u05dc\u05d5\u05e1\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05de\u05d4 \u05e2\u05e8\u05db\u05d9\u05dd \u05e9\u05dc k \u05dc- top-k \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 (\u05d1\u05de\u05e7\u05d5\u05e8 \u05d9\u05e9 \u05e2\u05e8\u05da k \u05d0\u05d7\u05d3). \u05db\u05db\u05d4 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05d1\u05db\u05de\u05d4 \u05d2\u05d3\u05dc\u05d9\u05dd \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 (\u05d7\u05d5\u05e5 \u05de\u05d6\u05d4 \u05d0\u05d9\u05df \u05d4\u05e8\u05d1\u05d4 \u05d3\u05de\u05d9\u05d5\u05df \u05db\u05d9 \u05d4\u05de\u05d8\u05e8\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d4\u05e4\u05d9\u05e7 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05d3\u05dc\u05d9\u05dc). \u05d4\u05d3\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e9\u05dc \u05d0\u05d9\u05d1\u05e8 \u05d4\u05de\u05e0\u05e1\u05d4 \u05dc\u05d2\u05e8\u05d5\u05dd \u05dc\u05e9\u05d2\u05d9\u05d0\u05ea \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e2\u05d1\u05d5\u05e8 top-k \u05e9\u05dc \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05de\u05ea\u05d9\u05dd (\u05e2\u05e8\u05db\u05d9\u05dd \u05d4\u05db\u05d9 \u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8) \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e9\u05d2\u05d9\u05d0\u05ea \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e9\u05dc \u05d4-top-k \u05e9\u05dc \u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d5\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8. \u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05e8\u05d3\u05ea \u05dc\u05e2\u05d5\u05de\u05e7 \u05d3\u05e2\u05ea\u05dd \u05dc\u05de\u05d4 \u05d6\u05d4 \u05e2\u05d5\u05d6\u05e8.\"\n$d.Paragraphs.Item(7).Range.Text = \"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05d0\u05d9\u05d1\u05e8 \u05e9\u05dc \u05d4\u05dc\u05d5\u05e1 \u05d4\u05e0\u05d9\u05d2\u05d5\u05d3\u05d9 \u05dc\u05d6\u05d4 \u05e9\u05de\u05ea\u05d5\u05d0\u05e8 \u05d1\u05e4\u05e1\u05e7\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea\u2026.\"\n$d.Paragraphs.Item(8).Range.Text = \"\u05d8\u05d5\u05d1, \u05e0\u05db\u05d5\u05df \u05e9\u05d4\u05d5\u05e4\u05d9\u05e2\u05d4 \u05dc\u05e0\u05d5 \u05d4\u05de\u05d8\u05e8\u05d9\u05d5\u05e9\u05e7\u05d4 \u05d1\u05e9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05d1\u05d9\u05e0\u05d5 \u05dc\u05d1\u05d9\u05df \u05d4\u05de\u05d8\u05e8\u05d9\u05e9\u05e7\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05d3\u05d9 \u05e8\u05d5\u05e4\u05e3. \u05d0\u05d1\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d7\u05d5\u05e5 \u05de\u05d6\u05d4\u2026.\"\n\n# --- Paragraphs 9-23 (1-based): delete entirely (removed body of the old review) ---\nfor ($i = 23; $i -ge 9; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# --- Paragraph 24 (1-based, now last paragraph, the link): replace URL text ---\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"https://arxiv.org/pdf/2503.01776\"\n"}
